$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 4: PyhtonHomePath value changes to a new local python install path
$ws.Range("B4").Value = "C:\Users\adelinas\AppData\Local\Programs\Python\Python38"

# Row 6: ScriptWorkingFolder value changes back to the TestingTool_v4 root folder
$ws.Range("B6").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4"

# Insert four new rows before the current row 10 (RobotModelFullPath) to make
# room for the two new script-path rows (DFSSymbolic / Concolic), each
# followed by its own spacer row, and push RobotModelFullPath further down.
$ws.Rows("10:13").Insert()
$ws.Rows("10:13").RowHeight = 14.25

# Row 8 becomes ScriptFullPathOfflineAll
$ws.Range("A8").Value = "ScriptFullPathOfflineAll"
$ws.Range("B8").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\bankLoan_offlineall.py"

# Row 10: ScriptFullPathDFSSymbolic (new)
$ws.Range("A10").Value = "ScriptFullPathDFSSymbolic"
$ws.Range("B10").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\bankLoan_dfssymbolic.py"

# Row 12: ScriptFullPathConcolic (new)
$ws.Range("A12").Value = "ScriptFullPathConcolic"
$ws.Range("B12").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\bankLoan_concolic.py"

# Row 13 is a pure spacer row (no styled-but-empty cells), matching the
# blank separator rows used throughout this sheet.
$ws.Range("A13:C13").Clear()

# Row 14 (was row 10): RobotModelFullPath stays the same text, already shifted by insert
$ws.Range("B14").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\Applications\C#Models\SimpleBankLoanCSharp"

$ws.Range("B12").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
